# Updated symbol list on Tue Dec 13 22:27:48 UTC 2022 with GitHub Actions
#
# This refreshes the cryptocurrency price snapshot: most rows only get a
# new value in column D (Price). A couple of rows also pick up new
# "Worstin24h"/"Bestin24h" suffixes in column E, and rows 49/50 swap places
# (BOLO <-> CryptobidCoin) together with their own refreshed prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText {
    param($Cell, $Text)
    # Column D holds numeric-looking strings (e.g. "0.08294", "0.00001400")
    # that must be preserved verbatim -- force text storage first so Excel
    # doesn't silently coerce them into floating point numbers.
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
}

# --- Column D (Price) refreshes -------------------------------------------
Set-PriceText $ws.Range("D2")  "270.96"
Set-PriceText $ws.Range("D4")  "6.338"
Set-PriceText $ws.Range("D6")  "3.652"
Set-PriceText $ws.Range("D7")  "6.692"
Set-PriceText $ws.Range("D8")  "1.387"
Set-PriceText $ws.Range("D9")  "0.8304"
Set-PriceText $ws.Range("D11") "0.1603"
Set-PriceText $ws.Range("D12") "0.08294"
Set-PriceText $ws.Range("D13") "0.03436"
Set-PriceText $ws.Range("D14") "0.03226"
Set-PriceText $ws.Range("D15") "0.09340"
Set-PriceText $ws.Range("D16") "3.850"
Set-PriceText $ws.Range("D17") "0.001652"
Set-PriceText $ws.Range("D18") "0.04730"
Set-PriceText $ws.Range("D19") "0.006304"
Set-PriceText $ws.Range("D20") "0.005668"
Set-PriceText $ws.Range("D21") "0.001077"
Set-PriceText $ws.Range("D23") "3.719"
Set-PriceText $ws.Range("D24") "2.399"
Set-PriceText $ws.Range("D25") "0.3347"
Set-PriceText $ws.Range("D27") "0.0002705"
Set-PriceText $ws.Range("D40") "0.04687"
Set-PriceText $ws.Range("D41") "0.007036"
Set-PriceText $ws.Range("D42") "0.1162"

# Row 43 (CEJI): price update + new "Worstin24h" tag appended to column E
Set-PriceText $ws.Range("D43") "0.003291"
$ws.Range("E43").Value = "42CEJICEJIWorstin24h"

Set-PriceText $ws.Range("D44") "0.01200"
Set-PriceText $ws.Range("D45") "0.00006276"
Set-PriceText $ws.Range("D46") "0.0009904"
Set-PriceText $ws.Range("D47") "0.00000000750"
Set-PriceText $ws.Range("D48") "0.9202"

# --- Rows 49/50 swap: BOLO and CryptobidCoin trade places ------------------
# Row 49 becomes CryptobidCoin (was BOLO)
$ws.Range("B49").Value = "CryptobidCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/h39bvStAP+cryptobidcoin-cbc"
Set-PriceText $ws.Range("D49") "0.00001400"
$ws.Range("E49").Value = "48CryptobidCoinCBC"

# Row 50 becomes BOLO (was CryptobidCoin)
$ws.Range("B50").Value = "BOLO"
$ws.Range("C50").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
Set-PriceText $ws.Range("D50") "0.001416"
$ws.Range("E50").Value = "49BOLOBOLO"

Set-PriceText $ws.Range("D51") "0.01240"
